$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F6").Value = "ppe"
$ws.Range("F7").Value = "ppe"
$ws.Range("F8").Value = "ppe"
$ws.Range("F12").Value = "env warning - water"
$ws.Range("F13").Value = "application instructions || env warning - species"
$ws.Range("F19").Value = "application instructions"
$ws.Range("F20").Value = "application instructions"
$ws.Range("F21").Value = "application instructions"
$ws.Range("F22").Value = "135_product_information"
$ws.Range("F23").Value = "use restrictions"
$ws.Range("F25").Value = "application instructions"
$ws.Range("F26").Value = "mixing"
$ws.Range("F27").Value = "mixing"
$ws.Range("F34").Value = "use restrictions"
$ws.Range("F38").Value = "off target movement"
$ws.Range("F39").Value = "off target movement"
$ws.Range("F41").Value = "use restrictions || off target movement"
$ws.Range("F42").Value = "application instructions"
$ws.Range("F43").Value = "application instructions"
$ws.Range("F44").Value = "application instructions"
$ws.Range("F45").Value = "application instructions"
$ws.Range("F46").Value = "application instructions"
$ws.Range("F47").Value = "mixing"
$ws.Range("F48").Value = "mixing"
$ws.Range("F50").Value = "mixing"
$ws.Range("F51").Value = "application instructions"
$ws.Range("F52").Value = "application instructions"
$ws.Range("F53").Value = "application instructions"
$ws.Range("F54").Value = "irrigation || application instructions || chemigation"
$ws.Range("F55").Value = "safety procedures"
$ws.Range("F56").Value = "safety procedures"
$ws.Range("F58").Value = "safety procedures"
$ws.Range("F60").Value = "irrigation"
$ws.Range("F61").Value = "irrigation"
$ws.Range("F62").Value = "irrigation"
$ws.Range("F64").Value = "irrigation"
$ws.Range("F332").Value = "154_pesticide_storage"
